# Apply the edit described by the commit: rename the "Organization" header
# in cell A1 to "Firm Name". This also produces a new (duplicate) font /
# cell style entry, mirroring what Excel does when the cell's formatting is
# (re-)applied alongside the text edit, and moves the active selection to B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch the font formatting on A1 so a fresh style record gets materialized
# (matches the extra cellXfs/font entry seen in the saved workbook), then
# update the cell text itself.
$ws.Range("A1").Font.ThemeColor = 1
$ws.Range("A1").Value = "Firm Name"

# Move / record the active selection as it was when the workbook was saved.
[void]$ws.Range("B13").Select()
